# While-Increment.pptx — "Added PDF versions of slides" edit
#
# 1) Title slide subtitle: "Spring 20XX" is split into two runs,
#    "Spring/Fall " and "20XX", so the term can later be edited
#    independently of the year.
#
# 2) The canonical diff also re-stamps the cached datetimeFigureOut
#    field text (handout master + notes master) from 6/17/2021 to
#    8/16/2021. Those are auto-computed "date last saved/printed"
#    placeholders that PowerPoint recalculates itself whenever the
#    deck is touched/printed; we don't hand-author them here.

$p = $ppt.ActivePresentation

# --- 1) Title slide: "Spring 20XX" -> "Spring/Fall " + "20XX" ------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) {
            continue
        }

        $tr = $shape.TextFrame.TextRange
        $fullText = $tr.Text
        $needle = "Spring 20XX"
        $pos = $fullText.IndexOf($needle)

        if ($pos -ge 0) {
            $prefix = "Spring "
            $start = $pos + 1          # TextRange.Characters is 1-based
            $len = $prefix.Length

            # Re-write just the "Spring " prefix so the run splits into
            # "Spring/Fall " + "20XX" (the "20XX" run is left untouched).
            $target = $tr.Characters($start, $len)
            $target.Text = "Spring/Fall "
        }
    }
}
